# Contest 16 RCB vs RR
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter this contest's raw points for each player (row 28)
$ws.Range("E28").Value = 0
$ws.Range("H28").Value = 70
$ws.Range("K28").Value = 50
$ws.Range("N28").Value = 80
$ws.Range("Q28").Value = 40
$ws.Range("T28").Value = 60
$ws.Range("W28").Value = 30
$ws.Range("Z28").Value = 100
$ws.Range("AC28").Value = 20

# Fix the Sundar (P) column's ranking formula for rows 28-42 so it matches the
# other players' formula pattern (use the full 9-player tie-break list and the
# A2:B10 lookup table instead of the legacy "score" name / 6-player list).
for ($r = 28; $r -le 42; $r++) {
    $formula = "=IF(ISERROR(VLOOKUP(RANK(Q$r, (`$AC$r,`$Z$r,`$W$r,`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$B`$10, 2, FALSE)),`"`",VLOOKUP(RANK(Q$r, (`$AC$r,`$Z$r,`$W$r,`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$B`$10, 2, FALSE))"
    $ws.Range("P$r").Formula = $formula
}

# Extend the season-total SUM formulas to cover through row 42
$ws.Range("E45").Formula = "=SUM(D13:D42)"
$ws.Range("H45").Formula = "=SUM(G13:G42)"
$ws.Range("K45").Formula = "=SUM(J13:J42)"
$ws.Range("N45").Formula = "=SUM(M13:M42)"
$ws.Range("Q45").Formula = "=SUM(P13:P42)"
$ws.Range("T45").Formula = "=SUM(S13:S42)"
$ws.Range("W45").Formula = "=SUM(V13:V42)"
$ws.Range("Z45").Formula = "=SUM(Y13:Y42)"
$ws.Range("AC45").Formula = "=SUM(AB13:AB42)"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("I45").Select()
$ws.Activate()
